$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet has 5 data rows (rows 2-5). The two oldest entries
# (rows 2 and 3) are being removed, and the remaining two entries
# (previously rows 4 and 5) shift up to become rows 2 and 3.
$ws.Rows("2:3").Delete()
